$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-17 15:42:58"
$wsZh.Range("G2").Value = "2016-01-17 15:43:42"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-17 15:43:09"
$wsDe.Range("G2").Value = "2016-01-17 15:44:01"
